$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing text storage so values that look
# numeric (e.g. "616.36", "1.00", "0.999") are not silently coerced into
# actual numbers (which would also lose precision / trailing zeros).
function Set-TextValue($range, $value) {
    $ws.Range($range).NumberFormat = "@"
    $ws.Range($range).Value = $value
}

# --- Simple value/percentage updates (rows that only change D and/or E) ---

# Row 2
Set-TextValue "D2" "69.154.05"
$ws.Range("E2").Value = "  +0.16%  "

# Row 3
Set-TextValue "D3" "3.734.11"
$ws.Range("E3").Value = "  -0.43%  "

# Row 4
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
Set-TextValue "D5" "616.36"
$ws.Range("E5").Value = "  +5.59%  "

# Row 6
Set-TextValue "D6" "186.30"
$ws.Range("E6").Value = "  +3.35%  "

# Row 7
Set-TextValue "D7" "3.730.68"
$ws.Range("E7").Value = "  -0.31%  "

# Row 8
$ws.Range("E8").Value = "  -0.85%  "

# Row 9
$ws.Range("E9").Value = "  -0.10%  "

# Row 10
$ws.Range("E10").Value = "  -1.19%  "

# Row 11
$ws.Range("E11").Value = "  -5.04%  "

# Row 12
Set-TextValue "D12" "56.68"
$ws.Range("E12").Value = "  +4.04%  "

# Row 13
$ws.Range("E13").Value = "  -5.05%  "

# Row 14
Set-TextValue "D14" "10.70"
$ws.Range("E14").Value = "  -2.03%  "

# Row 15
Set-TextValue "D15" "4.318.46"
$ws.Range("E15").Value = "  -0.58%  "

# Row 16
Set-TextValue "D16" "3.731.42"
$ws.Range("E16").Value = "  -1.41%  "

# Row 17
Set-TextValue "D17" "19.48"
$ws.Range("E17").Value = "  -1.48%  "

# Row 18
Set-TextValue "D18" "13.11"
$ws.Range("E18").Value = "  -1.60%  "

# Row 19
$ws.Range("E19").Value = "  -0.85%  "

# Row 20
$ws.Range("E20").Value = "  -2.03%  "

# Row 21
Set-TextValue "D21" "68.930.12"
$ws.Range("E21").Value = "  +0.17%  "

# Row 22
Set-TextValue "D22" "414.50"
$ws.Range("E22").Value = "  -0.84%  "

# Row 23
Set-TextValue "D23" "4.67"
$ws.Range("E23").Value = "  +0.70%  "

# Row 24
Set-TextValue "D24" "89.53"
$ws.Range("E24").Value = "  -1.32%  "

# Row 25
$ws.Range("E25").Value = "  -1.73%  "

# Row 26
Set-TextValue "D26" "12.85"
$ws.Range("E26").Value = "  -2.80%  "

# Row 27
Set-TextValue "D27" "10.97"
$ws.Range("E27").Value = "  -0.84%  "

# Row 28
$ws.Range("E28").Value = "  +2.14%  "

# Row 29
$ws.Range("E29").Value = "  -1.10%  "

# Row 30
Set-TextValue "D30" "9.69"
$ws.Range("E30").Value = "  -2.04%  "

# Row 31
$ws.Range("E31").Value = "  -1.20%  "

# Row 32
Set-TextValue "D32" "7.34"
$ws.Range("E32").Value = "  -13.66%  "

# Row 33
Set-TextValue "D33" "12.75"
$ws.Range("E33").Value = "  -2.29%  "

# Row 34
$ws.Range("E34").Value = "  +1.52%  "

# --- Row 35: Bittensor -> OKB (name/link/value swap) ---
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D35" "66.16"
$ws.Range("E35").Value = "  -1.03%  "

# Row 36
Set-TextValue "D36" "44.35"
$ws.Range("E36").Value = "  -3.64%  "

# --- Row 37: OKB -> Bittensor (name/link/value swap) ---
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D37" "613.91"
$ws.Range("E37").Value = "  +0.99%  "

# Row 38
Set-TextValue "D38" "0.0₃0870"
$ws.Range("E38").Value = "  -9.30%  "

# Row 39
Set-TextValue "D39" "0.408"
$ws.Range("E39").Value = "  -0.92%  "

# Row 40
Set-TextValue "D40" "1.00"
$ws.Range("E40").Value = "  +0.10%  "

# Row 41
Set-TextValue "D41" "0.999"
$ws.Range("E41").Value = "  +0.12%  "

# Row 42
$ws.Range("E42").Value = "  +0.82%  "

# Row 43
$ws.Range("E43").Value = "  -2.18%  "

# Row 44
$ws.Range("E44").Value = "  -1.33%  "

# Row 45
Set-TextValue "D45" "2.65"
$ws.Range("E45").Value = "  -1.17%  "

# Row 46
Set-TextValue "D46" "0.141"
$ws.Range("E46").Value = "  +1.98%  "

# --- Row 47: THORChain -> Maker (name/link/value swap) ---
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D47" "2.846.78"
$ws.Range("E47").Value = "  +2.37%  "

# --- Row 48: Maker -> THORChain (name/link/value swap) ---
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D48" "9.24"
$ws.Range("E48").Value = "  -4.50%  "

# Row 49
$ws.Range("E49").Value = "  +1.39%  "

# Row 50
$ws.Range("E50").Value = "  -17.98%  "

# Row 51
Set-TextValue "D51" "3.10"
$ws.Range("E51").Value = "  -3.28%  "
